$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date style (format) from A464 down to the new date cells A465:A491
$ws.Range("A464").Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)

$newRows = @(
  @(465, 44539, 4, 35, 310.0363185401719),
  @(466, 44540, 1, 31, 274.6035964212951),
  @(467, 44541, 2, 25, 221.4545132429799),
  @(468, 44542, 5, 27, 239.1708743024183),
  @(469, 44543, 8, 29, 256.8872353618567),
  @(470, 44544, 2, 24, 212.5963327132607),
  @(471, 44545, 1, 23, 203.7381521835415),
  @(472, 44546, 10, 29, 256.8872353618567),
  @(473, 44547, 8, 36, 318.894499069891),
  @(474, 44548, 3, 37, 327.7526795996102),
  @(475, 44550, 11, 43, 380.9017627779254),
  @(476, 44551, 5, 40, 354.3272211887678),
  @(477, 44552, 0, 38, 336.6108601293294),
  @(478, 44553, 15, 52, 460.6253875453982),
  @(479, 44554, 11, 53, 469.4835680751174),
  @(480, 44555, 13, 58, 513.7744707237133),
  @(481, 44556, 17, 72, 637.7889981397821),
  @(482, 44557, 5, 66, 584.639914961467),
  @(483, 44558, 1, 62, 549.2071928425902),
  @(484, 44559, 6, 68, 602.3562760209053),
  @(485, 44560, 28, 81, 717.5126229072549),
  @(486, 44561, 30, 100, 885.8180529719195),
  @(487, 44562, 24, 111, 983.2580387988307),
  @(488, 44563, 14, 108, 956.6834972096731),
  @(489, 44564, 16, 119, 1054.123483036584),
  @(490, 44565, 9, 127, 1124.988927274338),
  @(491, 44566, 15, 136, 1204.712552041811)
)

foreach ($row in $newRows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Host "Done. Last row dimension check:" $ws.UsedRange.Address()
